$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder "Polinesia Francesa" ahead of "Botsuana" / "Mali" (country list re-sorted by case count) ---
$ws.Range("A147").Value = "Polinesia Francesa"
$ws.Range("A148").Value = "Botsuana"
$ws.Range("A149").Value = "Mali"

# --- Update statistic values ---
$ws.Range("B4").Value = 8148532
$ws.Range("C4").Value = 58182
$ws.Range("D4").Value = 5269905
$ws.Range("E4").Value = 2656803
$ws.Range("G4").Value = 951
$ws.Range("H4").Value = 221824

$ws.Range("D5").Value = 6380456
$ws.Range("E5").Value = 813303

$ws.Range("B11").Value = 856951
$ws.Range("C11").Value = 2977
$ws.Range("D11").Value = 759597
$ws.Range("E11").Value = 63842
$ws.Range("G11").Value = 93
$ws.Range("H11").Value = 33512

$ws.Range("B30").Value = 189385
$ws.Range("C30").Value = 2504
$ws.Range("D30").Value = 159351
$ws.Range("E30").Value = 20370
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = 9664

$ws.Range("B67").Value = 51845
$ws.Range("C67").Value = 648
$ws.Range("D67").Value = 33921
$ws.Range("E67").Value = 16793
$ws.Range("G67").Value = 23
$ws.Range("H67").Value = 1131

$ws.Range("B117").Value = 7572
$ws.Range("C117").Value = 7
$ws.Range("D117").Value = 7314
$ws.Range("E117").Value = 95

$ws.Range("B132").Value = 5083
$ws.Range("C132").Value = 11
$ws.Range("D132").Value = 4886
$ws.Range("E132").Value = 89
$ws.Range("G132").Value = 1
$ws.Range("H132").Value = 108

$ws.Range("B147").Value = 3573
$ws.Range("C147").Value = 322
$ws.Range("D147").Value = 2487
$ws.Range("E147").Value = 1073
$ws.Range("G147").Value = 2
$ws.Range("H147").Value = 13

$ws.Range("B148").Value = 3515
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 853
$ws.Range("E148").Value = 2642
$ws.Range("H148").Value = 20

$ws.Range("B149").Value = 3352
$ws.Range("C149").Value = 55
$ws.Range("D149").Value = 2550
$ws.Range("E149").Value = 670
$ws.Range("H149").Value = 132

$ws.Range("B162").Value = 1983
$ws.Range("C162").Value = 11
$ws.Range("D162").Value = 1473
$ws.Range("E162").Value = 459
$ws.Range("G162").Value = 1
$ws.Range("H162").Value = 51

$ws.Range("D169").Value = 896
$ws.Range("E169").Value = 18

$ws.Range("B171").Value = 729
$ws.Range("C171").Value = 10
$ws.Range("D171").Value = 657
$ws.Range("E171").Value = 50

$ws.Range("D177").Value = 497
$ws.Range("E177").Value = 31

$ws.Range("B190").Value = 225
$ws.Range("C190").Value = 4
$ws.Range("D190").Value = 212
$ws.Range("E190").Value = 12

$ws.Range("B191").Value = 215
$ws.Range("C191").Value = 5
$ws.Range("D191").Value = 192

$ws.Range("B196").Value = 147
$ws.Range("C196").Value = 1
$ws.Range("E196").Value = 1

$ws.Range("B206").Value = 31
$ws.Range("C206").Value = 2
$ws.Range("E206").Value = 4

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Octubre de 2020 a las 02:36"
